$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.180.00'
$ws.Range('E2').Value = '  +6.25%  '
$ws.Range('D3').Value = '3.110.37'
$ws.Range('E3').Value = '  +4.08%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.55'
$ws.Range('E5').Value = '  +4.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.10'
$ws.Range('E6').Value = '  +3.96%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.103.05'
$ws.Range('E8').Value = '  +4.22%  '
$ws.Range('E9').Value = '  +2.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.150'
$ws.Range('E10').Value = '  +13.19%  '
$ws.Range('E11').Value = '  +8.17%  '
$ws.Range('E12').Value = '  +3.09%  '
$ws.Range('E13').Value = '  +7.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.51'
$ws.Range('E14').Value = '  +5.17%  '
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('D16').Value = '3.624.18'
$ws.Range('E16').Value = '  +4.08%  '
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '63.061.97'
$ws.Range('E18').Value = '  +6.07%  '
$ws.Range('D19').Value = '3.106.73'
$ws.Range('E19').Value = '  +3.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '466.36'
$ws.Range('E20').Value = '  +7.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.16'
$ws.Range('E21').Value = '  +4.32%  '
$ws.Range('E22').Value = '  +1.11%  '
$ws.Range('E23').Value = '  +7.41%  '
$ws.Range('E24').Value = '  -0.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.80'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.38'
$ws.Range('E27').Value = '  +8.21%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('E29').Value = '  +5.14%  '
$ws.Range('E31').Value = '  +9.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.96'
$ws.Range('E32').Value = '  +4.72%  '
$ws.Range('E33').Value = '  +3.11%  '
$ws.Range('D34').Value = '0.0₃0863'
$ws.Range('E34').Value = '  +11.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.43'
$ws.Range('E35').Value = '  +16.59%  '
$ws.Range('E36').Value = '  +4.98%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.32'
$ws.Range('E37').Value = '  +19.32%  '
$ws.Range('E38').Value = '  +3.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.86'
$ws.Range('E39').Value = '  +4.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '432.30'
$ws.Range('E40').Value = '  +8.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.72'
$ws.Range('E41').Value = '  +0.62%  '
$ws.Range('D42').Value = '2.919.48'
$ws.Range('E42').Value = '  +6.00%  '
$ws.Range('E43').Value = '  +4.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.279'
$ws.Range('E44').Value = '  +11.53%  '
$ws.Range('E45').Value = '  +6.01%  '
$ws.Range('E46').Value = '  +7.60%  '
$ws.Range('E47').Value = '  +1.95%  '
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '122.82'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.50'
$ws.Range('E51').Value = '  +4.74%  '
